$p = $ppt.ActivePresentation

# The decorative "💡 ..." caption text boxes are auto-fit-to-text shapes
# (<a:bodyPr wrap="none"><a:spAutoFit/></a:bodyPr>). Re-writing their
# TextRange recomputes a layout-engine height estimate, so after editing
# their text we restore the original geometry (identical on every slide
# that has one of these boxes: off=(5486400,3200400) ext=(2743200,914400)
# EMU == Left 432pt, Top 252pt, Width 216pt, Height 72pt).
function Reset-CaptionBoxGeometry($shape) {
    $shape.Left = 432
    $shape.Top = 252
    $shape.Width = 216
    $shape.Height = 72
}

# ---------------------------------------------------------------------------
# Slide 2 - "Key Takeaways / Summary": just update the three bullet strings.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Paragraphs(2).Text = ""
$tr2.Paragraphs(2).Text = "AI is transforming healthcare delivery and research."
$tr2.Paragraphs(3).Text = ""
$tr2.Paragraphs(3).Text = "AI applications range from diagnostics to drug discovery."
$tr2.Paragraphs(4).Text = ""
$tr2.Paragraphs(4).Text = "Ethical considerations are paramount for responsible AI adoption."

# ---------------------------------------------------------------------------
# Slide 3 - "Introduction to AI in Healthcare" -> "Introduction: AI Revolutionizing Healthcare"
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tt3 = $s3.Shapes.Item(1).TextFrame.TextRange
$tt3.Paragraphs(1).Text = ""
$tt3.Paragraphs(1).Text = "Introduction: AI Revolutionizing Healthcare"

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
$tr3.Paragraphs(2).Text = ""
$tr3.Paragraphs(2).Text = "AI is rapidly changing the healthcare landscape, offering unprecedented opportunities to improve patient care, streamline processes, and accelerate research."
$tr3.Paragraphs(3).Text = ""
$tr3.Paragraphs(3).Text = "Enhanced Diagnostics"
[void]$tr3.InsertAfter("`rPersonalized Treatment Plans`rDrug Discovery Acceleration`rOperational Efficiency")

$cap3 = $s3.Shapes.Item(3)
$cap3.TextFrame.TextRange.Paragraphs(1).Text = ""
$cap3.TextFrame.TextRange.Paragraphs(1).Text = "💡 AI brain scan visualization"
Reset-CaptionBoxGeometry $cap3

# ---------------------------------------------------------------------------
# Slide 4 - "AI-Powered Diagnostics" -> "AI in Diagnostics: Early and Accurate"
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tt4 = $s4.Shapes.Item(1).TextFrame.TextRange
$tt4.Paragraphs(1).Text = ""
$tt4.Paragraphs(1).Text = "AI in Diagnostics: Early and Accurate"

$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
$tr4.Paragraphs(2).Text = ""
$tr4.Paragraphs(2).Text = "AI algorithms excel at analyzing medical images (X-rays, MRIs, CT scans) to detect diseases earlier and more accurately."
$tr4.Paragraphs(3).Text = ""
$tr4.Paragraphs(3).Text = "Detecting Cancerous Tumors"
$tr4.Paragraphs(4).Text = ""
$tr4.Paragraphs(4).Text = "Identifying Anomalies in Imaging"
[void]$tr4.InsertAfter("`rReducing Diagnostic Errors`rImproving Patient Outcomes")

$cap4 = $s4.Shapes.Item(3)
$cap4.TextFrame.TextRange.Paragraphs(1).Text = ""
$cap4.TextFrame.TextRange.Paragraphs(1).Text = "💡 AI analyzing medical image"
Reset-CaptionBoxGeometry $cap4

# ---------------------------------------------------------------------------
# Slide 5 - "Drug Discovery and Development" -> "Personalized Medicine: Tailoring Treatment"
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$tt5 = $s5.Shapes.Item(1).TextFrame.TextRange
$tt5.Paragraphs(1).Text = ""
$tt5.Paragraphs(1).Text = "Personalized Medicine: Tailoring Treatment"

$tr5 = $s5.Shapes.Item(2).TextFrame.TextRange
$tr5.Paragraphs(2).Text = ""
$tr5.Paragraphs(2).Text = "AI analyzes patient data (genetics, lifestyle, medical history) to create personalized treatment plans."
$tr5.Paragraphs(3).Text = ""
$tr5.Paragraphs(3).Text = "Predicting Treatment Response"
$tr5.Paragraphs(4).Text = ""
$tr5.Paragraphs(4).Text = "Optimizing Drug Dosage"
[void]$tr5.InsertAfter("`rIdentifying Patients at Risk`rImproving Treatment Efficacy")

$cap5 = $s5.Shapes.Item(3)
$cap5.TextFrame.TextRange.Paragraphs(1).Text = ""
$cap5.TextFrame.TextRange.Paragraphs(1).Text = "💡 DNA helix personalized treatment"
Reset-CaptionBoxGeometry $cap5

# ---------------------------------------------------------------------------
# Slide 6 - "Personalized Medicine with AI" -> "Drug Discovery: Accelerating Innovation"
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$tt6 = $s6.Shapes.Item(1).TextFrame.TextRange
$tt6.Paragraphs(1).Text = ""
$tt6.Paragraphs(1).Text = "Drug Discovery: Accelerating Innovation"

$tr6 = $s6.Shapes.Item(2).TextFrame.TextRange
$tr6.Paragraphs(2).Text = ""
$tr6.Paragraphs(2).Text = "AI accelerates drug discovery by analyzing vast datasets to identify potential drug candidates and predict their efficacy."
$tr6.Paragraphs(3).Text = ""
$tr6.Paragraphs(3).Text = "Target Identification"
$tr6.Paragraphs(4).Text = ""
$tr6.Paragraphs(4).Text = "Drug Design"
[void]$tr6.InsertAfter("`rClinical Trial Optimization`rReducing Development Time")

$cap6 = $s6.Shapes.Item(3)
$cap6.TextFrame.TextRange.Paragraphs(1).Text = ""
$cap6.TextFrame.TextRange.Paragraphs(1).Text = "💡 AI designing drug molecule"
Reset-CaptionBoxGeometry $cap6

# ---------------------------------------------------------------------------
# Slide 7 - "AI in Robotic Surgery" -> "Operational Efficiency: Streamlining Processes"
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$tt7 = $s7.Shapes.Item(1).TextFrame.TextRange
$tt7.Paragraphs(1).Text = ""
$tt7.Paragraphs(1).Text = "Operational Efficiency: Streamlining Processes"

$tr7 = $s7.Shapes.Item(2).TextFrame.TextRange
$tr7.Paragraphs(2).Text = ""
$tr7.Paragraphs(2).Text = "AI automates administrative tasks, optimizes workflows, and improves resource allocation in healthcare settings."
$tr7.Paragraphs(3).Text = ""
$tr7.Paragraphs(3).Text = "Automating Scheduling"
$tr7.Paragraphs(4).Text = ""
$tr7.Paragraphs(4).Text = "Predicting Patient Volume"
[void]$tr7.InsertAfter("`rOptimizing Bed Management`rReducing Healthcare Costs")

$cap7 = $s7.Shapes.Item(3)
$cap7.TextFrame.TextRange.Paragraphs(1).Text = ""
$cap7.TextFrame.TextRange.Paragraphs(1).Text = "💡 Hospital automation dashboard display"
Reset-CaptionBoxGeometry $cap7

# ---------------------------------------------------------------------------
# Slide 8 - "Challenges and Ethical Considerations" -> "Ethical Considerations: Responsible AI"
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tt8 = $s8.Shapes.Item(1).TextFrame.TextRange
$tt8.Paragraphs(1).Text = ""
$tt8.Paragraphs(1).Text = "Ethical Considerations: Responsible AI"

$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange
$tr8.Paragraphs(2).Text = ""
$tr8.Paragraphs(2).Text = "Ethical considerations are crucial for responsible AI adoption in healthcare."
$tr8.Paragraphs(3).Text = ""
$tr8.Paragraphs(3).Text = "Data Privacy and Security"
$tr8.Paragraphs(4).Text = ""
$tr8.Paragraphs(4).Text = "Algorithmic Bias"
[void]$tr8.InsertAfter("`rTransparency and Explainability`rHuman Oversight and Accountability")

$cap8 = $s8.Shapes.Item(3)
$cap8.TextFrame.TextRange.Paragraphs(1).Text = ""
$cap8.TextFrame.TextRange.Paragraphs(1).Text = "💡 AI ethics privacy security"
Reset-CaptionBoxGeometry $cap8

# ---------------------------------------------------------------------------
# Slide 9 - "The Future of AI in Healthcare" -> "Challenges and Future Directions"
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$tt9 = $s9.Shapes.Item(1).TextFrame.TextRange
$tt9.Paragraphs(1).Text = ""
$tt9.Paragraphs(1).Text = "Challenges and Future Directions"

$tr9 = $s9.Shapes.Item(2).TextFrame.TextRange
$tr9.Paragraphs(2).Text = ""
$tr9.Paragraphs(2).Text = "Data availability and quality"
$tr9.Paragraphs(3).Text = ""
$tr9.Paragraphs(3).Text = "Regulatory hurdles"
$tr9.Paragraphs(4).Text = ""
$tr9.Paragraphs(4).Text = "Integration with existing systems"
[void]$tr9.InsertAfter("`rBuilding trust and acceptance among healthcare professionals and patients`rFocus on explainable AI (XAI)`rContinuous monitoring and evaluation")

$cap9 = $s9.Shapes.Item(3)
$cap9.TextFrame.TextRange.Paragraphs(1).Text = ""
$cap9.TextFrame.TextRange.Paragraphs(1).Text = "💡 Future healthcare technology vision"
Reset-CaptionBoxGeometry $cap9
